$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: row data (Fecha, Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg o Unidades) gets reshuffled across rows 2-32 as
# part of the daily/weekly logic consolidation. Apply the new values cell by
# cell.

$ws.Range("D2").Value = 44391
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 21000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21500
$ws.Range("P2").Value = 538
$ws.Range("D3").Value = 44742
$ws.Range("J3").Value = 120
$ws.Range("D4").Value = 44483
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 362
$ws.Range("D5").Value = 44167
$ws.Range("H5").Value = 'Española'
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13500
$ws.Range("N5").Value = '$/caja 30 unidades'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 450
$ws.Range("Q5").Value = 30
$ws.Range("D6").Value = 44377
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 21000
$ws.Range("M6").Value = 20333
$ws.Range("O6").Value = 'Región de Coquimbo'
$ws.Range("P6").Value = 508
$ws.Range("D7").Value = 44377
$ws.Range("H7").Value = 'Symphony'
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 21000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 21500
$ws.Range("P7").Value = 538
$ws.Range("D8").Value = 44859
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("O8").Value = 'Provincia de Limarí'
$ws.Range("P8").Value = 388
$ws.Range("D9").Value = 44398
$ws.Range("J9").Value = 170
$ws.Range("D10").Value = 44433
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("N10").Value = '$/caja 40 unidades'
$ws.Range("P10").Value = 488
$ws.Range("Q10").Value = 40
$ws.Range("D11").Value = 44419
$ws.Range("H11").Value = 'Symphony'
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 21000
$ws.Range("L11").Value = 22000
$ws.Range("M11").Value = 21500
$ws.Range("N11").Value = '$/caja 50 unidades'
$ws.Range("P11").Value = 430
$ws.Range("Q11").Value = 50
$ws.Range("D12").Value = 44356
$ws.Range("H12").Value = 'Argentina(o)'
$ws.Range("J12").Value = 120
$ws.Range("N12").Value = '$/caja 50 unidades'
$ws.Range("P12").Value = 390
$ws.Range("Q12").Value = 50
$ws.Range("D13").Value = 44435
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19500
$ws.Range("P13").Value = 488
$ws.Range("D14").Value = 44489
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("N14").Value = '$/caja 40 unidades'
$ws.Range("P14").Value = 338
$ws.Range("Q14").Value = 40
$ws.Range("D15").Value = 44827
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 362
$ws.Range("D16").Value = 44405
$ws.Range("J16").Value = 200
$ws.Range("D17").Value = 44160
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("I17").Value = 'Primera'
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("P17").Value = 362
$ws.Range("Q17").Value = 40
$ws.Range("D18").Value = 44370
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 20000
$ws.Range("L18").Value = 21000
$ws.Range("M18").Value = 20429
$ws.Range("N18").Value = '$/caja 50 unidades'
$ws.Range("O18").Value = 'Región de Coquimbo'
$ws.Range("P18").Value = 409
$ws.Range("Q18").Value = 50
$ws.Range("D19").Value = 44370
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 22000
$ws.Range("L19").Value = 23000
$ws.Range("M19").Value = 22500
$ws.Range("P19").Value = 562
$ws.Range("D20").Value = 44384
$ws.Range("H20").Value = 'Madrigal'
$ws.Range("J20").Value = 80
$ws.Range("D21").Value = 44384
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19333
$ws.Range("N21").Value = '$/caja 50 unidades'
$ws.Range("P21").Value = 387
$ws.Range("Q21").Value = 50
$ws.Range("D22").Value = 44384
$ws.Range("H22").Value = 'Symphony'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 21000
$ws.Range("M22").Value = 20400
$ws.Range("P22").Value = 510
$ws.Range("D23").Value = 44482
$ws.Range("H23").Value = 'Madrigal'
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("N23").Value = '$/caja 40 unidades'
$ws.Range("O23").Value = 'Región de Coquimbo'
$ws.Range("P23").Value = 362
$ws.Range("Q23").Value = 40
$ws.Range("D24").Value = 44769
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17500
$ws.Range("P24").Value = 438
$ws.Range("D25").Value = 44785
$ws.Range("H25").Value = 'Argentina(o)'
$ws.Range("I25").Value = 'Segunda'
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15500
$ws.Range("N25").Value = '$/caja 50 unidades'
$ws.Range("P25").Value = 310
$ws.Range("Q25").Value = 50
$ws.Range("D26").Value = 44363
$ws.Range("H26").Value = 'Madrigal'
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 19000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 19500
$ws.Range("N26").Value = '$/caja 40 unidades'
$ws.Range("P26").Value = 488
$ws.Range("Q26").Value = 40
$ws.Range("D27").Value = 44426
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 19000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 19500
$ws.Range("P27").Value = 488
$ws.Range("D28").Value = 44706
$ws.Range("J28").Value = 250
$ws.Range("D29").Value = 44468
$ws.Range("H29").Value = 'Argentina(o)'
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 17500
$ws.Range("N29").Value = '$/caja 50 unidades'
$ws.Range("P29").Value = 350
$ws.Range("Q29").Value = 50
$ws.Range("D30").Value = 44762
$ws.Range("H30").Value = 'Madrigal'
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 19000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 19500
$ws.Range("N30").Value = '$/caja 40 unidades'
$ws.Range("P30").Value = 488
$ws.Range("Q30").Value = 40
$ws.Range("D31").Value = 44806
$ws.Range("H31").Value = 'Argentina(o)'
$ws.Range("J31").Value = 250
$ws.Range("K31").Value = 14000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 14500
$ws.Range("O31").Value = 'Provincia de Limarí'
$ws.Range("P31").Value = 362
$ws.Range("D32").Value = 44412
$ws.Range("H32").Value = 'Symphony'
$ws.Range("J32").Value = 240
$ws.Range("K32").Value = 21000
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = 21500
$ws.Range("N32").Value = '$/caja 40 unidades'
$ws.Range("P32").Value = 538
$ws.Range("Q32").Value = 40
